$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("fasta-method-1")

# --- Rows 2-4 (num procs group "23" -> "1"): drop the per-step timing/memory
#     columns C:F, keeping just the num-procs/step label columns ---
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Index ref fasta"
$ws.Range("C2:F2").ClearContents()

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "Calcs"
$ws.Range("C3:F3").ClearContents()

$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "Add unmapped"
$ws.Range("C4:F4").ClearContents()

# --- Row 5 (num procs group "1", unchanged) also loses its C:F data ---
$ws.Range("A5").Value = 1
$ws.Range("B5").Value = "End program"
$ws.Range("C5:F5").ClearContents()

# --- New group: num procs = 4, rows 7-10 (row 6 left blank as a separator) ---
$ws.Range("A7").Value = 4
$ws.Range("B7").Value = "Index ref fasta"

$ws.Range("A8").Value = 4
$ws.Range("B8").Value = "Calcs"

$ws.Range("A9").Value = 4
$ws.Range("B9").Value = "Add unmapped"

$ws.Range("A10").Value = 4
$ws.Range("B10").Value = "End program"

# --- New group: num procs = 8, rows 12-15 (row 11 left blank as a separator) ---
$ws.Range("A12").Value = 8
$ws.Range("B12").Value = "Index ref fasta"
$ws.Range("C12").Value = 155.26003408400001
$ws.Range("D12").Value = 155.273874044
$ws.Range("E12").Value = 38.97265625
$ws.Range("F12").Value = 168.95703125

$ws.Range("A13").Value = 8
$ws.Range("B13").Value = "Calcs"
$ws.Range("C13").Value = 17.8446338177
$ws.Range("D13").Value = 173.118507862
$ws.Range("E13").Value = 49.5078125
$ws.Range("F13").Value = 179.46484375

$ws.Range("A14").Value = 8
$ws.Range("B14").Value = "Add unmapped"
$ws.Range("C14").Value = 0.0019459724426299999
$ws.Range("D14").Value = 6648.2537269599998
$ws.Range("E14").Value = 893.34765625
$ws.Range("F14").Value = 2287.35546875

$ws.Range("A15").Value = 8
$ws.Range("B15").Value = "End program"
$ws.Range("C15").Value = 0.0041451454162599998
$ws.Range("D15").Value = 9053.6434659999995
$ws.Range("E15").Value = 972
$ws.Range("F15").Value = 2366.125

# --- Group num procs = 23 moves down to rows 17-20 (row 16 left blank) ---
$ws.Range("A17").Value = 23
$ws.Range("B17").Value = "Index ref fasta"
$ws.Range("C17").Value = 23.881839036900001
$ws.Range("D17").Value = 23.890377044699999
$ws.Range("E17").Value = 39.5
$ws.Range("F17").Value = 184.46484375

$ws.Range("A18").Value = 23
$ws.Range("B18").Value = "Calcs"
$ws.Range("C18").Value = 15.988245964100001
$ws.Range("D18").Value = 43.822856903100003
$ws.Range("E18").Value = 49.09765625
$ws.Range("F18").Value = 178.96484375

$ws.Range("A19").Value = 23
$ws.Range("B19").Value = "Add unmapped"
$ws.Range("C19").Value = 0.00362300872803
$ws.Range("D19").Value = 2820.4306268700002
$ws.Range("E19").Value = 2298.38671875
$ws.Range("F19").Value = 5693.8984375

$ws.Range("A20").Value = 23
$ws.Range("B20").Value = "End program"
$ws.Range("C20").Value = 0.0062611103057900004
$ws.Range("D20").Value = 4913.9100429999999
$ws.Range("E20").Value = 2446.16015625
$ws.Range("F20").Value = 5841.78125

# --- Selection moves to G6 on this sheet ---
$ws.Range("G6").Select()
